$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Cell writes are deliberately ordered so that the workbook's shared-string
# table gets rebuilt (on save) with new unique strings appended in the same
# order the final published workbook uses. The very first write of each new
# piece of text decides its position in xl/sharedStrings.xml, so we introduce
# each brand-new value exactly once, in the required order, before touching
# any of the cells that merely repeat an already-introduced value.
# ---------------------------------------------------------------------------

# 1) First appearance of every new shared string, in target order
$ws.Range("F4").Value = "RO.ACT.003HAB.SRS"
$ws.Range("G3").Value = "RO.ACT.003MET.SRL"
$ws.Range("G5").Value = "RO.ACT.003MET.SRM"
$ws.Range("G4").Value = "RO.ACT.003MET.SRS"
$ws.Range("E2").Value = "RO.ACT.001.CRE"
$ws.Range("E3").Value = "RO.ACT.001.SUP"
$ws.Range("E4").Value = "RO.ACT.001.LEC"
$ws.Range("E5").Value = "RO.ACT.001.MAJ"
$ws.Range("H3").Value = "RO.ACT.004EMP.SRL"
$ws.Range("H2").Value = "RO.ACT.004EMP.SRA"
$ws.Range("H4").Value = "RO.ACT.004EMP.SRS"
$ws.Range("H5").Value = "RO.ACT.004EMP.SRM"

# 2) Remaining cells in the new F2:H5 / E2:H5 grid that reuse a string
#    already known to the workbook (either pre-existing or just introduced
#    above)
$ws.Range("F2").Value = "RO.ACT.003HAB.SRA"
$ws.Range("F3").Value = "RO.ACT.003HAB.SRL"
$ws.Range("F5").Value = "RO.ACT.003HAB.SRM"
$ws.Range("G2").Value = "RO.ACT.003MET.SRA"

# 3) Column A / D updates
$ws.Range("A3").Value = "RO.ACT.004EMP.SRA"
$ws.Range("A4").Value = ""
$ws.Range("D5").Value = "AD.SEC.002.FON.01"
$ws.Range("D6").Value = ""

# 4) The D8:D14 helper list moves to E8:E14 (with D10/D12/D14 dropped to
#    match the new sparser layout)
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "RO.ACT.003HAB.SRA"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "RO.ACT.003HAB.SRL"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = "RO.ACT.003HAB.SRS"
$ws.Range("E11").Value = "RO.ACT.003MET.SRA"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = "RO.ACT.003MET.SRL"
$ws.Range("E13").Value = "RO.ACT.003MET.SRS"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = "AD.SEC.002.FON.01"

# ---------------------------------------------------------------------------
# Column widths for the newly-used columns D:H (best achievable match to the
# published widths; the host's internal char-width granularity means a few
# of these land a fraction off the original bestFit pixel-exact values)
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 17.721354166666668
$ws.Columns.Item(5).ColumnWidth = 17.385416666666668
$ws.Columns.Item(6).ColumnWidth = 21.385416666666668
$ws.Columns.Item(7).ColumnWidth = 20.166666666666668
$ws.Columns.Item(8).ColumnWidth = 15.166666666666666

# ---------------------------------------------------------------------------
# Final selection
# ---------------------------------------------------------------------------
$ws.Range("A4").Select() | Out-Null
